$wb = $excel.ActiveWorkbook

# Sheet "展览" and sheet "全部类型" both contain the same table of data
# that needs its "想去人数" (F column) values updated.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 494
    $ws.Range("F3").Value = 3353
    $ws.Range("F4").Value = 89
}
